$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Statistik")

# Insert a new column before the current column G (pt_no), shifting it to H
$ws.Columns("G:G").Insert()

# Copy the header cell style from the old header (now H1, originally G1) to the new G1
$ws.Range("H1").Copy()
$ws.Range("G1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# Set new header and value
$ws.Range("G1").Value = "icdRd_no_ext"
$ws.Range("G2").Value = 297

$wb.Save()
